$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D136").Value = 0.791995474
$ws.Range("D137").Value = 0.788120887
$ws.Range("D138").Value = 0.597740902
$ws.Range("D139").Value = 0.620527487
$ws.Range("C140").Value = 0.241887844
$ws.Range("C141").Value = 0.331651578
$ws.Range("C142").Value = 0.154182215
$ws.Range("C143").Value = 0.166899468
$ws.Range("C144").Value = 0.042359665
$ws.Range("C145").Value = 0.266698307
